$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.906
$ws.Range("B9").Value = 5.312
$ws.Range("C9").Value = -10.623
$ws.Range("B18").Value = 5.275
$ws.Range("B20").Value = 7.069
$ws.Range("C23").Value = -12.953
$ws.Range("C24").Value = -12.548
$ws.Range("C26").Value = -12.807
$ws.Range("B27").Value = 5.752000000000001
$ws.Range("C34").Value = -12.196
$ws.Range("C35").Value = -12.316
$ws.Range("C48").Value = -11.861
$ws.Range("C52").Value = -11.743
$ws.Range("C66").Value = -11.574
$ws.Range("C67").Value = -10.875
$ws.Range("B69").Value = 5.827
$ws.Range("B76").Value = 6.247999999999999
$ws.Range("C80").Value = -12.534
$ws.Range("B82").Value = 5.486
$ws.Range("C99").Value = -12.048
